$d = $word.ActiveDocument

# --- Change 1 ---------------------------------------------------------
# "...and you need a reason why doesn't want to cooperate."
#   -> "...and you need a reason why he doesn't want to cooperate."
$r1 = $d.Content
$found1 = $r1.Find.Execute(
    "a reason why doesn" + [char]0x2019 + "t want to cooperate",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "a reason why he doesn" + [char]0x2019 + "t want to cooperate",
    2)

# --- Change 2 ---------------------------------------------------------
# Remove " And nothing the players can do will affect that." from the
# "Not having a reason to help..." paragraph.
$r2 = $d.Content
$found2 = $r2.Find.Execute(
    " And nothing the players can do will affect that.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 2)

# --- Change 3 ---------------------------------------------------------
# "...something to be said about being [pagebreak]evocative while not
# being overly verbose, but still giving..."
#   -> "...something to be said about being evocative while not being
#       overly verbose, but [pagebreak]still giving..."
# Move the phrase "evocative while not being overly verbose, but " from
# the start of the run that follows the page break up into the run that
# precedes it.
$r3 = $d.Content
$found3 = $r3.Find.Execute(
    "something to be said about being ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0)
$insertPoint = $r3.End
$phrase = "evocative while not being overly verbose, but "
$r3.Collapse(0)
$r3.InsertBefore($phrase)

$dupStart = $insertPoint + $phrase.Length
$dupEnd = $dupStart + $phrase.Length
$dupRange = $d.Range($dupStart, $dupEnd)
$dupRange.Delete()

Write-Output "change1=$found1 change2=$found2 change3=$found3"
